$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the title in C1 (English sheet title), dropping the period after "4.3.1.1"
$ws.Range("C1").Value = "4.3.1.1 Youth education by gender"

# 2. Add the new "2021" column (M) and extend column L's formatting to match the
#    rest of the data row (copy format from column K, the preceding data column),
#    then fill in the actual values for L (existing 2020 data, unchanged values)
#    and M (new 2021 data).

# Row 2 (separator row, no values, just formatting)
$ws.Range("K2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

# Row 3 (year headers)
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2020
$ws.Range("M3").Value = 2021

# Row 4
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 10.8

# Row 5
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("L5").Value = 6.4
$ws.Range("M5").Value = 5.2

# Row 6
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("L6").Value = 13.5
$ws.Range("M6").Value = 16.2

# Row 7
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("L7").Value = 24.3
$ws.Range("M7").Value = 24.2

# Row 8
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("L8").Value = 27.8
$ws.Range("M8").Value = 27.6

# Row 9
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("L9").Value = 20.9
$ws.Range("M9").Value = 20.9

# Row 10
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("L10").Value = 26.7
$ws.Range("M10").Value = 28.5

# Row 11
$ws.Range("K11").Copy()
$ws.Range("L11").PasteSpecial(-4122)
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("L11").Value = 28.4
$ws.Range("M11").Value = 29.7

# Row 12
$ws.Range("K12").Copy()
$ws.Range("L12").PasteSpecial(-4122)
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("L12").Value = 25
$ws.Range("M12").Value = 27.5

# 3. Restore the selection/active cell as it was left in the saved workbook
$ws.Range("O2").Select()
